$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet to reflect the new "through" date
$ws.Name = "Through 2022-09-17"

# Update the "September (through ...)" label in column A, row 10
$ws.Range("A10").Value = "September (through 09-17)"

# Update September row (row 10) counts for each year column B..I
$ws.Range("B10").Value = 17
$ws.Range("C10").Value = 30
$ws.Range("D10").Value = 40
$ws.Range("E10").Value = 32
$ws.Range("F10").Value = 39
$ws.Range("G10").Value = 63
$ws.Range("H10").Value = 87
$ws.Range("I10").Value = 82

# Update Total row (row 11) counts for each year column B..I
$ws.Range("B11").Value = 211
$ws.Range("C11").Value = 411
$ws.Range("D11").Value = 591
$ws.Range("E11").Value = 522
$ws.Range("F11").Value = 388
$ws.Range("G11").Value = 847
$ws.Range("H11").Value = 1157
$ws.Range("I11").Value = 1217
